$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D is overwritten with the "SVR" predictions that used to live in C
# (Random-Forest predictions now occupy column E instead).
$ws.Range("D2").Value = 3.13
$ws.Range("E2").Value = 3.06
$ws.Range("D3").Value = 3.35
$ws.Range("E3").Value = 3.17
$ws.Range("D4").Value = 3.22
$ws.Range("E4").Value = 3.16
$ws.Range("D5").Value = 3.31
$ws.Range("E5").Value = 3.19
$ws.Range("D6").Value = 3.2
$ws.Range("E6").Value = 3.21
$ws.Range("D7").Value = 3.32
$ws.Range("E7").Value = 3.09
$ws.Range("D8").Value = 3.2
$ws.Range("E8").Value = 3.28
$ws.Range("D9").Value = 3.07
$ws.Range("E9").Value = 2.99
$ws.Range("D10").Value = 3.29
$ws.Range("E10").Value = 2.96
$ws.Range("D11").Value = 3.19
$ws.Range("E11").Value = 2.89
$ws.Range("D12").Value = 3.19
$ws.Range("E12").Value = 2.93
$ws.Range("D13").Value = 3.08
$ws.Range("E13").Value = 3.06
$ws.Range("D14").Value = 3.08
$ws.Range("E14").Value = 2.97
$ws.Range("D15").Value = 3.21
$ws.Range("E15").Value = 2.96
$ws.Range("D16").Value = 3.23
$ws.Range("E16").Value = 2.99
$ws.Range("D17").Value = 3.11
$ws.Range("E17").Value = 2.92
$ws.Range("D18").Value = 3.05
$ws.Range("E18").Value = 3.02
$ws.Range("D19").Value = 3.05
$ws.Range("E19").Value = 3.08
$ws.Range("D20").Value = 1.9
$ws.Range("E20").Value = 1.7
$ws.Range("D21").Value = 2.33
$ws.Range("E21").Value = 2.66
$ws.Range("D22").Value = 2.34
$ws.Range("E22").Value = 2.7
$ws.Range("D23").Value = 2.3
$ws.Range("E23").Value = 2.72
$ws.Range("D24").Value = 2.3
$ws.Range("E24").Value = 2.72
$ws.Range("D25").Value = 2.36
$ws.Range("E25").Value = 2.72
$ws.Range("D26").Value = 7.01
$ws.Range("E26").Value = 6.41
$ws.Range("D27").Value = 3.62
$ws.Range("E27").Value = 3.77
$ws.Range("D28").Value = 3.75
$ws.Range("E28").Value = 2.99
$ws.Range("D29").Value = 2.7
$ws.Range("E29").Value = 2.76
$ws.Range("E30").Value = 3.76
$ws.Range("D31").Value = 8.51
$ws.Range("E31").Value = 8.18
$ws.Range("D32").Value = 2.54
$ws.Range("E32").Value = 2.71
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 2.65
$ws.Range("D34").Value = 3.04
$ws.Range("E34").Value = 3.11
$ws.Range("D35").Value = 3.67
$ws.Range("E35").Value = 3.11

# The active selection ends up parked on column D (user clicked there).
$ws.Columns("D").Select()
